{"js": "// Highlight quantitative impact metrics (percentages, dollar amounts, large\n// numbers) in bold + a dark slate color (#2C3E50) across the resume body.\n//\n// For each target paragraph we search for the specific metric substrings\n// (in left-to-right order of appearance) and apply bold + color formatting\n// to just those substrings, leaving the surrounding text in separate,\n// unformatted runs \u2014 mirroring the OOXML run-splitting seen in the diff.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Map of paragraph index (0-based, matching context.document.body.paragraphs)\n// -> ordered list of metric substrings that should become bold + colored.\nconst targets = [\n  { index: 9, metrics: [\"23%\", \"64%\"] },\n  { index: 11, metrics: [\"\u00b14.2%\", \"\u00b12.1%\", \"71%\", \"87%\"] },\n  { index: 12, metrics: [\"73.5%\", \"$4.7M\"] },\n  { index: 13, metrics: [\"$2\"] },\n  { index: 49, metrics: [\"\u00b14.2%\", \"\u00b12.1%\"] },\n  { index: 50, metrics: [\"71%\", \"87%\"] },\n  { index: 52, metrics: [\"34%\", \"28%\"] },\n];\n\nconst HIGHLIGHT_COLOR = \"#2C3E50\";\n\n// Kick off all the searches first, then sync once, then format, then sync\n// again \u2014 avoids re-searching against text whose run layout has already\n// changed by an earlier mutation within the same paragraph.\nconst searchResults = [];\nfor (const { index, metrics } of targets) {\n  const paragraph = paragraphs.items[index];\n  for (const metric of metrics) {\n    const found = paragraph.search(metric, { matchCase: true });\n    found.load(\"text\");\n    searchResults.push(found);\n  }\n}\n\nawait context.sync();\n\nfor (const found of searchResults) {\n  found.items[0].font.set({ bold: true, color: HIGHLIGHT_COLOR });\n}\n\nawait context.sync();\n", "ps1": "# Highlight quantitative impact metrics (percentages, dollar amounts, large\n# numbers) in bold + a dark slate color (#2C3E50) across the resume body.\n#\n# For each target paragraph we Find each metric substring (left to right)\n# and apply bold + color formatting to just that substring, leaving the\n# surrounding text as separate, unformatted runs.\n\nfunction RGB($r, $g, $b) {\n    return $r + ($g * 256) + ($b * 65536)\n}\n$highlightColor = RGB 0x2C 0x3E 0x50\n\n$d = $word.ActiveDocument\n\nfunction Highlight-Metric($paragraph, $metricText) {\n    $rng = $paragraph.Range\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $metricText\n    $find.MatchCase = $true\n    $find.Forward = $true\n    $find.Wrap = 0\n    $result = $find.Execute()\n    if ($result) {\n        $rng.Font.Bold = $true\n        $rng.Font.Color = $highlightColor\n    }\n}\n\n# Paragraph indices are 1-based (Word COM convention) and correspond to the\n# same paragraphs touched by the diff:\n#   10 -> \"Discovered systematic race coding errors ... from 23% to 64%\"\n#   12 -> \"Utilized advanced sampling methods ... \u00b14.2% to \u00b12.1% ... 71% to 87% ...\"\n#   13 -> \"Trigonometric algorithm ... 73.5% ... $4.7M ...\"\n#   14 -> \"Built real-time FEC analysis systems ... valued over $2 trillion\"\n#   50 -> \"Predictive excellence: ... \u00b14.2% to \u00b12.1%\"\n#   51 -> \"Increased voter turnout prediction accuracy from 71% to 87%\"\n#   53 -> \"Methodological advancement: Improved segmentation accuracy 34% and survey incidence 28%\"\n$targets = @(\n    @{ Index = 10; Metrics = @('23%', '64%') },\n    @{ Index = 12; Metrics = @('\u00b14.2%', '\u00b12.1%', '71%', '87%') },\n    @{ Index = 13; Metrics = @('73.5%', '$4.7M') },\n    @{ Index = 14; Metrics = @('$2') },\n    @{ Index = 50; Metrics = @('\u00b14.2%', '\u00b12.1%') },\n    @{ Index = 51; Metrics = @('71%', '87%') },\n    @{ Index = 53; Metrics = @('34%', '28%') }\n)\n\nforeach ($target in $targets) {\n    $paragraph = $d.Paragraphs.Item($target.Index)\n    foreach ($metric in $target.Metrics) {\n        Highlight-Metric $paragraph $metric\n    }\n}\n"}
